$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (rows 2-8) from serial date 45212 (2023-10-13) to 45221 (2023-10-22)
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45221
}
